$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 246
$ws.Range("I2").Value = 633
$ws.Range("J2").Value = 2531
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 714
$ws.Range("M2").Value = 46
$ws.Range("N2").Value = 439
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 39
$ws.Range("S2").Value = 296
$ws.Range("T2").Value = 455
$ws.Range("U2").Value = 33
$ws.Range("V2").Value = 3986
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 3947
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 61
$ws.Range("AA2").Value = 28
